$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the shared "LAST SCRAPE DATE" value for every data row (F2:F81).
# This mirrors the single shared-string edit "2019-03-07" -> "2019-03-12".
# Force the range to text first so Excel keeps it as a string instead of
# auto-converting it to a date serial number, then drop the format override
# again so no stray per-cell style survives in the saved file.
$ws.Range("F2:F81").NumberFormat = "@"
$ws.Range("F2:F81").Value = "2019-03-12"
$ws.Range("F2:F81").NumberFormat = "General"
$ws.Range("F2:F81").ClearFormats()

# Re-order several GAME NAME / GAME NUMBER / TOP PRIZES REMAINING rows
# (some new games were interleaved into the price-tier blocks, shifting
# the rows below them) and update the one genuine TOP PRIZES REMAINING
# count that changed (row 67, $50 & $100 BLOWOUT: 19340 -> 17844).
$ws.Cells.Item(9, 3).Value = 'LOOSE CHANGE'
$ws.Cells.Item(9, 4).Value = 237
$ws.Cells.Item(9, 5).Value = 7
$ws.Cells.Item(10, 3).Value = 'MISTLE DOUGH'
$ws.Cells.Item(10, 4).Value = 224
$ws.Cells.Item(10, 5).Value = 10
$ws.Cells.Item(11, 3).Value = 'ELECTRIC 8S'
$ws.Cells.Item(11, 4).Value = 226
$ws.Cells.Item(11, 5).Value = 22
$ws.Cells.Item(12, 3).Value = 'CRAZY 7s'
$ws.Cells.Item(12, 4).Value = 181
$ws.Cells.Item(12, 5).Value = 26
$ws.Cells.Item(14, 3).Value = 'SAPPHIRE 7S'
$ws.Cells.Item(14, 4).Value = 236
$ws.Cells.Item(14, 5).Value = 6
$ws.Cells.Item(15, 3).Value = '$100 IN A FLASH'
$ws.Cells.Item(15, 4).Value = 185
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(26, 3).Value = '10X THE MONEY'
$ws.Cells.Item(26, 4).Value = 173
$ws.Cells.Item(26, 5).Value = 0
$ws.Cells.Item(27, 3).Value = 'WIN ALL!'
$ws.Cells.Item(27, 4).Value = 168
$ws.Cells.Item(27, 5).Value = 0
$ws.Cells.Item(33, 3).Value = 'DIAMOND BINGO'
$ws.Cells.Item(33, 4).Value = 217
$ws.Cells.Item(33, 5).Value = 2
$ws.Cells.Item(34, 3).Value = 'BLOCK-O'
$ws.Cells.Item(34, 4).Value = 231
$ws.Cells.Item(34, 5).Value = 5
$ws.Cells.Item(45, 3).Value = 'MONEY CLIP'
$ws.Cells.Item(45, 4).Value = 233
$ws.Cells.Item(45, 5).Value = 3
$ws.Cells.Item(47, 3).Value = 'NATIONAL LAMPOONS CHRISTMAS VACATION (TM)'
$ws.Cells.Item(47, 4).Value = 221
$ws.Cells.Item(47, 5).Value = 2
$ws.Cells.Item(48, 3).Value = '$100,000 CASINO CASH'
$ws.Cells.Item(48, 4).Value = 197
$ws.Cells.Item(48, 5).Value = 1
$ws.Cells.Item(49, 3).Value = 'JUMBO CASH'
$ws.Cells.Item(49, 4).Value = 161
$ws.Cells.Item(49, 5).Value = 1
$ws.Cells.Item(50, 3).Value = '20X LUCKY'
$ws.Cells.Item(50, 4).Value = 228
$ws.Cells.Item(50, 5).Value = 4
$ws.Cells.Item(51, 3).Value = 'I LOVE LUCY'
$ws.Cells.Item(51, 4).Value = 235
$ws.Cells.Item(51, 5).Value = 3
$ws.Cells.Item(52, 3).Value = 'TRIPLE BONUS CROSSWORD'
$ws.Cells.Item(52, 4).Value = 246
$ws.Cells.Item(52, 5).Value = 5
$ws.Cells.Item(56, 3).Value = 'CASH CRAZE'
$ws.Cells.Item(56, 4).Value = 206
$ws.Cells.Item(56, 5).Value = 0
$ws.Cells.Item(57, 3).Value = 'DOUBLE DIAMOND'
$ws.Cells.Item(57, 4).Value = 189
$ws.Cells.Item(57, 5).Value = 0
$ws.Cells.Item(67, 3).Value = '$50 & $100 BLOWOUT'
$ws.Cells.Item(67, 4).Value = 196
$ws.Cells.Item(67, 5).Value = 17844
